# Extend the pension-statistics table with a new "2022" column (S), mirroring
# the existing 2021 column (R) formatting, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of column R (rows 2-6) into the new column S
# so the new cells inherit the same cell styles (borders, number formats,
# fonts, etc.) as the rest of the table, instead of falling back to the
# sheet's default column style.
$ws.Range("R2:R6").Copy()
$ws.Range("S2:S6").PasteSpecial(-4122)

# Populate the new 2022 data column.
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 265803
$ws.Range("S5").Value = 3.8
$ws.Range("S6").Value = 33.6

# Move the active cell/selection as recorded in the saved view state.
$ws.Range("C19").Select()
